$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H33").Value = 483.2069
$ws_ALC.Range("I33").Value = 340.63635
$ws_ALC.Range("K33").Value = 340.63635
$ws_ALC.Range("M33").Value = -111.63635

$ws_ALC.Range("H129").Value = 3907307.5
$ws_ALC.Range("J129").Value = 1027.7413
$ws_ALC.Range("L129").Value = 3083.2239
$ws_ALC.Range("N129").Value = -13083.2239

$ws_ALC.Range("H130").Value = 39333.332
$ws_ALC.Range("J130").Value = 39333.332
$ws_ALC.Range("L130").Value = 39333.332
$ws_ALC.Range("N130").Value = -49373.332

$ws_ALC.Range("H131").Value = 5799.5557
$ws_ALC.Range("I131").Value = 5257.5
$ws_ALC.Range("J131").Value = 6883.6665
$ws_ALC.Range("K131").Value = 15772.5
$ws_ALC.Range("L131").Value = 20650.9995
$ws_ALC.Range("M131").Value = -10732.5
$ws_ALC.Range("N131").Value = -30730.9995

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H5").Value = 50
$ws_ARM.Range("I5").Value = 50
$ws_ARM.Range("J5").Value = 50
$ws_ARM.Range("K5").Value = 50
$ws_ARM.Range("L5").Value = 50
$ws_ARM.Range("M5").Value = 62
$ws_ARM.Range("N5").Value = -274

$ws_ARM.Range("H61").Value = 2609.476
$ws_ARM.Range("I61").Value = 1011.92
$ws_ARM.Range("J61").Value = 4958.8237
$ws_ARM.Range("K61").Value = 1011.92
$ws_ARM.Range("L61").Value = 4958.8237
$ws_ARM.Range("M61").Value = -799.92
$ws_ARM.Range("N61").Value = -5382.8237

$ws_ARM.Range("H97").Value = 667.4761999999999
$ws_ARM.Range("I97").Value = 643
$ws_ARM.Range("J97").Value = 900
$ws_ARM.Range("K97").Value = 643
$ws_ARM.Range("L97").Value = 900
$ws_ARM.Range("M97").Value = -147
$ws_ARM.Range("N97").Value = -1892

$ws_ARM.Range("H102").Value = 2659
$ws_ARM.Range("I102").Value = 1848.75
$ws_ARM.Range("J102").Value = 5900
$ws_ARM.Range("K102").Value = 1848.75
$ws_ARM.Range("L102").Value = 5900
$ws_ARM.Range("M102").Value = -226.75
$ws_ARM.Range("N102").Value = -9144

$ws_ARM.Range("H110").Value = 1296.1333
$ws_ARM.Range("I110").Value = 669.38464
$ws_ARM.Range("J110").Value = 5370
$ws_ARM.Range("K110").Value = 669.38464
$ws_ARM.Range("L110").Value = 5370
$ws_ARM.Range("M110").Value = 1375.61536
$ws_ARM.Range("N110").Value = -9460

$ws_ARM.Range("H132").Value = 25002768
$ws_ARM.Range("I132").Value = 55557824
$ws_ARM.Range("J132").Value = 3177.5454
$ws_ARM.Range("K132").Value = 166673472
$ws_ARM.Range("L132").Value = 9532.636200000001
$ws_ARM.Range("M132").Value = -166670942
$ws_ARM.Range("N132").Value = -14592.6362

$ws_ARM.Range("H136").Value = 2609.476
$ws_ARM.Range("I136").Value = 1011.92
$ws_ARM.Range("J136").Value = 4958.8237
$ws_ARM.Range("K136").Value = 3035.76
$ws_ARM.Range("L136").Value = 14876.4711
$ws_ARM.Range("M136").Value = -485.7599999999998
$ws_ARM.Range("N136").Value = -19976.4711

$ws_ARM.Range("H137").Value = 0
$ws_ARM.Range("J137").Value = 0
$ws_ARM.Range("L137").Value = 0
$ws_ARM.Range("N137").ClearContents()

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H4").Value = 50
$ws_BSM.Range("I4").Value = 50
$ws_BSM.Range("J4").Value = 50
$ws_BSM.Range("K4").Value = 50
$ws_BSM.Range("L4").Value = 50
$ws_BSM.Range("M4").Value = 65
$ws_BSM.Range("N4").Value = -280

$ws_BSM.Range("H22").Value = 271.2857
$ws_BSM.Range("I22").Value = 99.75
$ws_BSM.Range("J22").Value = 500
$ws_BSM.Range("K22").Value = 99.75
$ws_BSM.Range("L22").Value = 500
$ws_BSM.Range("M22").Value = 73.25
$ws_BSM.Range("N22").Value = -846

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H132").Value = 1764.9822
$ws_CRP.Range("I132").Value = 1294.4
$ws_CRP.Range("J132").Value = 3690.0908
$ws_CRP.Range("K132").Value = 3883.2
$ws_CRP.Range("L132").Value = 11070.2724
$ws_CRP.Range("M132").Value = -1353.2
$ws_CRP.Range("N132").Value = -16130.2724

$ws_CRP.Range("H134").Value = 1195.3405
$ws_CRP.Range("I134").Value = 551.5952
$ws_CRP.Range("J134").Value = 6602.8
$ws_CRP.Range("K134").Value = 1654.7856
$ws_CRP.Range("L134").Value = 19808.4
$ws_CRP.Range("M134").Value = 880.2144000000001
$ws_CRP.Range("N134").Value = -24878.4

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H12").Value = 130.13637
$ws_CUL.Range("I12").Value = 9.777778
$ws_CUL.Range("J12").Value = 213.46153
$ws_CUL.Range("K12").Value = 29.333334
$ws_CUL.Range("L12").Value = 640.38459
$ws_CUL.Range("M12").Value = 143.666666
$ws_CUL.Range("N12").Value = -986.38459

$ws_CUL.Range("H138").Value = 2984.25
$ws_CUL.Range("J138").Value = 4124.2856
$ws_CUL.Range("L138").Value = 12372.8568
$ws_CUL.Range("N138").Value = -22652.8568

$ws_CUL.Range("H140").Value = 5750152
$ws_CUL.Range("I140").Value = 13889580
$ws_CUL.Range("J140").Value = 4673.5293
$ws_CUL.Range("K140").Value = 41668740
$ws_CUL.Range("L140").Value = 14020.5879
$ws_CUL.Range("M140").Value = -41663560
$ws_CUL.Range("N140").Value = -24380.5879

$ws_CUL.Range("H141").Value = 3194.1667
$ws_CUL.Range("I141").Value = 2592.2222
$ws_CUL.Range("J141").Value = 5000
$ws_CUL.Range("K141").Value = 7776.6666
$ws_CUL.Range("L141").Value = 15000
$ws_CUL.Range("M141").Value = -2596.6666
$ws_CUL.Range("N141").Value = -25360

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H2").Value = 32.705883
$ws_GSM.Range("I2").Value = 26.333334
$ws_GSM.Range("J2").Value = 48
$ws_GSM.Range("K2").Value = 26.333334
$ws_GSM.Range("L2").Value = 48
$ws_GSM.Range("M2").Value = 86.66666599999999
$ws_GSM.Range("N2").Value = -274

$ws_GSM.Range("H97").Value = 2568.158
$ws_GSM.Range("I97").Value = 1656.8462
$ws_GSM.Range("J97").Value = 4542.6665
$ws_GSM.Range("K97").Value = 1656.8462
$ws_GSM.Range("L97").Value = 4542.6665
$ws_GSM.Range("M97").Value = -1160.8462
$ws_GSM.Range("N97").Value = -5534.6665

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H93").Value = 3223.8462
$ws_LTW.Range("I93").Value = 2455.7778
$ws_LTW.Range("J93").Value = 4952
$ws_LTW.Range("K93").Value = 2455.7778
$ws_LTW.Range("L93").Value = 4952
$ws_LTW.Range("M93").Value = -1207.7778
$ws_LTW.Range("N93").Value = -7448

$ws_LTW.Range("H132").Value = 2194.83
$ws_LTW.Range("I132").Value = 1376.4138
$ws_LTW.Range("J132").Value = 3183.75
$ws_LTW.Range("K132").Value = 4129.2414
$ws_LTW.Range("L132").Value = 9551.25
$ws_LTW.Range("M132").Value = -1599.2414
$ws_LTW.Range("N132").Value = -14611.25

$ws_LTW.Range("H133").Value = 28189.572
$ws_LTW.Range("J133").Value = 28189.572
$ws_LTW.Range("L133").Value = 28189.572
$ws_LTW.Range("N133").Value = -33249.572

$ws_LTW.Range("H135").Value = 29734.875
$ws_LTW.Range("J135").Value = 29734.875
$ws_LTW.Range("L135").Value = 29734.875
$ws_LTW.Range("N135").Value = -39874.875

$ws_LTW.Range("H137").Value = 29486.363
$ws_LTW.Range("J137").Value = 29486.363
$ws_LTW.Range("L137").Value = 29486.363
$ws_LTW.Range("N137").Value = -39686.363

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H107").Value = 949.58826
$ws_WVR.Range("I107").Value = 437.5
$ws_WVR.Range("J107").Value = 1404.7778
$ws_WVR.Range("K107").Value = 1312.5
$ws_WVR.Range("L107").Value = 4214.3334
$ws_WVR.Range("M107").Value = 607.5
$ws_WVR.Range("N107").Value = -8054.3334

$ws_WVR.Range("H122").Value = 419203.5
$ws_WVR.Range("I122").Value = 557589.75
$ws_WVR.Range("J122").Value = 4044.8333
$ws_WVR.Range("K122").Value = 1672769.25
$ws_WVR.Range("L122").Value = 12134.4999
$ws_WVR.Range("M122").Value = -1670319.25
$ws_WVR.Range("N122").Value = -17034.4999

$ws_WVR.Range("H135").Value = 76738.336
$ws_WVR.Range("J135").Value = 76738.336
$ws_WVR.Range("L135").Value = 76738.336
$ws_WVR.Range("N135").Value = -86878.336

$ws_WVR.Range("H136").Value = 993.7451
$ws_WVR.Range("I136").Value = 646.55554
$ws_WVR.Range("J136").Value = 1827
$ws_WVR.Range("K136").Value = 1939.66662
$ws_WVR.Range("L136").Value = 5481
$ws_WVR.Range("M136").Value = 610.33338
$ws_WVR.Range("N136").Value = -10581
